$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analysis ver2")

# --- Reset the previous range (values + formatting) before rebuilding ---
$ws.Range("A1:X11").ClearContents()
$ws.Range("A1:X11").ClearFormats()

# --- Row 1: headers ---
$ws.Range("A1").Value = "Crow_A.wav"
$ws.Range("B1").Value = "Elephant_D.wav"
$ws.Range("C1").Value = "CarStart_E.wav"
$ws.Range("D1").Value = "Chomp_E.wav"
$ws.Range("E1").Value = "Coin_E.wav"
$ws.Range("F1").Value = "Clap_C.wav"
$ws.Range("G1").Value = "Cough_E.wav"
$ws.Range("H1").Value = "Initial/Repeat/Phone_A.wav"
$ws.Range("I1").Value = "Faucet_C.wav"
$ws.Range("J1").Value = "Chick_A.wav"
$ws.Range("K1").Value = "Dog_F.wav"
$ws.Range("L1").Value = "Duck_B.wav"
$ws.Range("M1").Value = "Goat_C.wav"
$ws.Range("N1").Value = "Growl_F.wav"
$ws.Range("O1").Value = "Laugh_D.wav"
$ws.Range("P1").Value = "Pour_C.wav"
$ws.Range("Q1").Value = "Siren_B.wav"
$ws.Range("R1").Value = "Whistle_C.wav"
$ws.Range("S1").Value = "Baby_E.wav"
$ws.Range("T1").Value = "Initial/Doubled/Chime_A.wav"
$ws.Range("U1").Value = "Snore_C.wav"
$ws.Range("V1").Value = "Test/Lures/Bubbles_B.wav"
$ws.Range("W1").Value = "Writing_C.wav"
$ws.Range("X1").Value = "Faucet_B.wav"
$ws.Range("Y1").Value = "Footsteps_C.wav"
$ws.Range("Z1").Value = "Footsteps_F.wav"

# --- Rows 2-6: round 1 data ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 2
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 3
$ws.Range("Y2").Value = 2
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 3
$ws.Range("X3").Value = 3
$ws.Range("Y3").Value = 3
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 2
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 3
$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 3
$ws.Range("Y4").Value = 2
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = 2
$ws.Range("X5").Value = 3
$ws.Range("Y5").Value = 2
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 2
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 2
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 2
$ws.Range("X6").Value = 2
$ws.Range("Y6").Value = 3

# --- Row 7 intentionally blank ---

# --- Rows 8-10: round 2 data ---
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 2
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 2
$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 2
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 2
$ws.Range("W8").Value = 2
$ws.Range("Z8").Value = 2
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 2
$ws.Range("M9").Value = 3
$ws.Range("N9").Value = 2
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 2
$ws.Range("U9").Value = 3
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = 1
$ws.Range("Z9").Value = 2
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 3
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 2
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 2
$ws.Range("P10").Value = 2
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = 2
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 2
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 3
$ws.Range("W10").Value = 2
$ws.Range("Z10").Value = 3

# --- Row 11 intentionally blank ---

# --- Row 12: proportion of rating "1" (most similar) ---
$ws.Range("A12").Formula = "=COUNTIF(A2:A10,1)/8"
$ws.Range("B12").Formula = "=COUNTIF(B2:B10,1)/8"
$ws.Range("C12").Formula = "=COUNTIF(C2:C6,1)/5"
$ws.Range("D12").Formula = "=COUNTIF(D2:D10,1)/8"
$ws.Range("E12").Formula = "=COUNTIF(E2:E10,1)/8"
$ws.Range("F12").Formula = "=COUNTIF(F2:F10,1)/8"
$ws.Range("G12").Formula = "=COUNTIF(G2:G10,1)/8"
$ws.Range("H12").Formula = "=COUNTIF(H2:H10,1)/8"
$ws.Range("I12").Formula = "=COUNTIF(I2:I10,1)/3"
$ws.Range("J12").Formula = "=COUNTIF(J2:J10,1)/8"
$ws.Range("K12").Formula = "=COUNTIF(K2:K10,1)/8"
$ws.Range("L12").Formula = "=COUNTIF(L2:L10,1)/8"
$ws.Range("M12").Formula = "=COUNTIF(M2:M10,1)/8"
$ws.Range("N12").Formula = "=COUNTIF(N2:N10,1)/8"
$ws.Range("O12").Formula = "=COUNTIF(O2:O10,1)/8"
$ws.Range("P12").Formula = "=COUNTIF(P2:P10,1)/8"
$ws.Range("Q12").Formula = "=COUNTIF(Q2:Q10,1)/8"
$ws.Range("R12").Formula = "=COUNTIF(R2:R10,1)/8"
$ws.Range("S12").Formula = "=COUNTIF(S2:S6,1)/5"
$ws.Range("T12").Formula = "=COUNTIF(T2:T10,1)/8"
$ws.Range("U12").Formula = "=COUNTIF(U2:U10,1)/8"
$ws.Range("V12").Formula = "=COUNTIF(V2:V10,1)/8"
$ws.Range("W12").Formula = "=COUNTIF(W2:W10,1)/8"
$ws.Range("X12").Formula = "=COUNTIF(X2:X10,1)/8"
$ws.Range("Y12").Formula = "=COUNTIF(Y2:Y10,1)/8"
$ws.Range("Z12").Formula = "=COUNTIF(Z2:Z10,1)/3"

# --- Row 13: proportion of rating "3" (least similar) ---
$ws.Range("A13").Formula = "=COUNTIF(A2:A10,3)/8"
$ws.Range("B13").Formula = "=COUNTIF(B2:B10,3)/8"
$ws.Range("C13").Formula = "=COUNTIF(C2:C6,3)/5"
$ws.Range("D13").Formula = "=COUNTIF(D2:D10,3)/8"
$ws.Range("E13").Formula = "=COUNTIF(E2:E10,3)/8"
$ws.Range("F13").Formula = "=COUNTIF(F2:F10,3)/8"
$ws.Range("G13").Formula = "=COUNTIF(G2:G10,3)/8"
$ws.Range("H13").Formula = "=COUNTIF(H2:H10,3)/8"
$ws.Range("I13").Formula = "=COUNTIF(I2:I10,3)/3"
$ws.Range("J13").Formula = "=COUNTIF(J2:J10,3)/8"
$ws.Range("K13").Formula = "=COUNTIF(K2:K10,3)/8"
$ws.Range("L13").Formula = "=COUNTIF(L2:L10,3)/8"
$ws.Range("M13").Formula = "=COUNTIF(M2:M10,3)/8"
$ws.Range("N13").Formula = "=COUNTIF(N2:N10,3)/8"
$ws.Range("O13").Formula = "=COUNTIF(O2:O10,3)/8"
$ws.Range("P13").Formula = "=COUNTIF(P2:P10,3)/8"
$ws.Range("Q13").Formula = "=COUNTIF(Q2:Q10,3)/8"
$ws.Range("R13").Formula = "=COUNTIF(R2:R10,3)/8"
$ws.Range("S13").Formula = "=COUNTIF(S2:S6,3)/5"
$ws.Range("T13").Formula = "=COUNTIF(T2:T10,3)/8"
$ws.Range("U13").Formula = "=COUNTIF(U2:U10,3)/8"
$ws.Range("V13").Formula = "=COUNTIF(V2:V10,3)/8"
$ws.Range("W13").Formula = "=COUNTIF(W2:W10,3)/8"
$ws.Range("X13").Formula = "=COUNTIF(X2:X10,3)/8"
$ws.Range("Y13").Formula = "=COUNTIF(Y2:Y10,3)/8"
$ws.Range("Z13").Formula = "=COUNTIF(Z2:Z10,3)/3"

# --- Row 14: proportion of rating "2" ---
$ws.Range("A14").Formula = "=COUNTIF(A2:A10,2)/8"
$ws.Range("B14").Formula = "=COUNTIF(B2:B10,2)/8"
$ws.Range("C14").Formula = "=COUNTIF(C2:C10,2)/8"
$ws.Range("D14").Formula = "=COUNTIF(D2:D10,2)/8"
$ws.Range("E14").Formula = "=COUNTIF(E2:E10,2)/8"
$ws.Range("F14").Formula = "=COUNTIF(F2:F10,2)/8"
$ws.Range("G14").Formula = "=COUNTIF(G2:G10,2)/8"
$ws.Range("H14").Formula = "=COUNTIF(H2:H10,2)/8"
$ws.Range("I14").Formula = "=COUNTIF(I2:I10,2)/8"
$ws.Range("J14").Formula = "=COUNTIF(J2:J10,2)/8"
$ws.Range("K14").Formula = "=COUNTIF(K2:K10,2)/8"
$ws.Range("L14").Formula = "=COUNTIF(L2:L10,2)/8"
$ws.Range("M14").Formula = "=COUNTIF(M2:M10,2)/8"
$ws.Range("N14").Formula = "=COUNTIF(N2:N10,2)/8"
$ws.Range("O14").Formula = "=COUNTIF(O2:O10,2)/8"
$ws.Range("P14").Formula = "=COUNTIF(P2:P10,2)/8"
$ws.Range("Q14").Formula = "=COUNTIF(Q2:Q10,2)/8"
$ws.Range("R14").Formula = "=COUNTIF(R2:R10,2)/8"
$ws.Range("S14").Formula = "=COUNTIF(S2:S10,2)/8"
$ws.Range("T14").Formula = "=COUNTIF(T2:T10,2)/8"
$ws.Range("U14").Formula = "=COUNTIF(U2:U10,2)/8"
$ws.Range("V14").Formula = "=COUNTIF(V2:V10,2)/8"
$ws.Range("W14").Formula = "=COUNTIF(W2:W10,2)/8"
$ws.Range("X14").Formula = "=COUNTIF(X2:X10,2)/8"
$ws.Range("Y14").Formula = "=COUNTIF(Y2:Y10,2)/8"
$ws.Range("Z14").Formula = "=COUNTIF(Z2:Z10,2)/8"

# --- Row 15: a few highlighted (yellow) blank marker cells ---
$ws.Range("D15").Interior.ColorIndex = 6
$ws.Range("I15").Interior.ColorIndex = 6
$ws.Range("N15").Interior.ColorIndex = 6
$ws.Range("O15").Interior.ColorIndex = 6
$ws.Range("P15").Interior.ColorIndex = 6
$ws.Range("U15").Interior.ColorIndex = 6
$ws.Range("Z15").Interior.ColorIndex = 6

# --- Column M (13th) width ---
$ws.Columns.Item(13).ColumnWidth = 8.7265625

# --- Clear the stale highlight fill from the old N/S/T/U/W columns ---
$ws.Range("A1:Z14").Interior.ColorIndex = -4142

# --- Selection / active cell like the source file ---
$ws.Range("G15").Select() | Out-Null
